{"js": "// Release notes proofing pass: split several paragraphs' runs at word\n// boundaries and add the w:proofErr spell-check / grammar-check markers\n// Word's proofing pass leaves behind, without changing the visible text\n// (the only actual formatting change is that \"testJCHashGet()\" loses the\n// Consolas/shaded run formatting while \" - \" keeps it).\n\nconst W_NS = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\n\nfunction pkg(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"' + W_NS + '\"><w:body>' +\n    bodyInnerXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst paragraphXml = {\n  0:\n    '<w:p><w:r><w:t xml:space=\"preserve\">Release notes for </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>JumpCloud</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Assignment</w:t></w:r></w:p>',\n\n  2:\n    '<w:p><w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>JCManual</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>TestCases.xlsx</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">  -</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> created for manual tests and requirements tracking</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>',\n\n  3:\n    \"<w:p><w:r><w:t>JumpCloudTest</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">.jar is the jar file for </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>api</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> automation on JAVA using </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>Eclipse(</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r><w:t>IDE)</w:t></w:r></w:p>\",\n\n  4:\n    '<w:p><w:r><w:t xml:space=\"preserve\">Automation script </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">\\u2013 </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>testJCHashGet</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r><w:t>(</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r><w:t>)</w:t></w:r>\" +\n    \"<w:r><w:rPr>\" +\n    '<w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\" w:cs=\"Consolas\"/>' +\n    '<w:color w:val=\"000000\"/>' +\n    '<w:sz w:val=\"32\"/>' +\n    '<w:szCs w:val=\"32\"/>' +\n    '<w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"E8F2FE\"/>' +\n    '</w:rPr><w:t xml:space=\"preserve\"> - </w:t></w:r></w:p>',\n\n  6:\n    '<w:p><w:r><w:t xml:space=\"preserve\">Password validation is not happening which is a </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>know</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> issue</w:t></w:r>' +\n    \"<w:r><w:t>.</w:t></w:r></w:p>\",\n\n  7:\n    \"<w:p><w:r><w:t>Issue of average time showing 0 yet time is observed and reported as part manual tests and failed the test along with the steps</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\"> and also included postman log for the </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>same</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/></w:p>',\n\n  9:\n    '<w:p><w:r><w:t xml:space=\"preserve\">Postman execution scripts are also included for the </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r><w:t>same</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/></w:p>',\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst indices = Object.keys(paragraphXml)\n  .map(Number)\n  .sort((a, b) => a - b);\n\nfor (const idx of indices) {\n  const para = paragraphs.items[idx];\n  para.insertOoxml(pkg(paragraphXml[idx]), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Release notes proofing pass: split several paragraphs' runs at word\n# boundaries and add the w:proofErr spell-check / grammar-check markers\n# Word's proofing pass leaves behind, without changing the visible text\n# (the only actual formatting change is that \"testJCHashGet()\" loses the\n# Consolas/shaded run formatting while \" - \" keeps it).\n\n$d = $word.ActiveDocument\n\nfunction New-ParaXml([string]$innerXml) {\n    return @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>$innerXml</w:body></w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n}\n\n$paragraphXml = @{\n    1 = '<w:p><w:r><w:t xml:space=\"preserve\">Release notes for </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>JumpCloud</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> Assignment</w:t></w:r></w:p>'\n\n    3 = '<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>JCManual</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>TestCases.xlsx</w:t></w:r><w:r><w:t xml:space=\"preserve\">  -</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> created for manual tests and requirements tracking</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r></w:p>'\n\n    4 = '<w:p><w:r><w:t>JumpCloudTest</w:t></w:r><w:r><w:t xml:space=\"preserve\">.jar is the jar file for </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>api</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> automation on JAVA using </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>Eclipse(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>IDE)</w:t></w:r></w:p>'\n\n    5 = '<w:p><w:r><w:t xml:space=\"preserve\">Automation script </w:t></w:r><w:r><w:t xml:space=\"preserve\">\u2013 </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>testJCHashGet</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\" w:cs=\"Consolas\"/><w:color w:val=\"000000\"/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"E8F2FE\"/></w:rPr><w:t xml:space=\"preserve\"> - </w:t></w:r></w:p>'\n\n    7 = '<w:p><w:r><w:t xml:space=\"preserve\">Password validation is not happening which is a </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>know</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> issue</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'\n\n    8 = '<w:p><w:r><w:t>Issue of average time showing 0 yet time is observed and reported as part manual tests and failed the test along with the steps</w:t></w:r><w:r><w:t xml:space=\"preserve\"> and also included postman log for the </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>same</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p>'\n\n    10 = '<w:p><w:r><w:t xml:space=\"preserve\">Postman execution scripts are also included for the </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>same</w:t></w:r><w:proofErr w:type=\"gramEnd\"/></w:p>'\n}\n\nforeach ($idx in ($paragraphXml.Keys | Sort-Object)) {\n    $lastIndex = $d.Paragraphs.Count\n    if ($idx -eq $lastIndex) {\n        # Word keeps a sentinel paragraph mark at the very end of the story;\n        # replacing the last paragraph's own Range with a full <w:p> (which\n        # carries its own mark) leaves a stray empty paragraph behind. Work\n        # around it: insert a new paragraph ahead of the last one, fill that\n        # (now not-last) paragraph with the target XML, then delete the old\n        # trailing paragraph outright.\n        $old = $d.Paragraphs($idx)\n        $old.Range.InsertParagraphBefore()\n        $d.Paragraphs($idx).Range.InsertXML((New-ParaXml $paragraphXml[$idx]))\n        $d.Paragraphs($idx + 1).Range.Delete()\n    } else {\n        $range = $d.Paragraphs($idx).Range\n        $range.InsertXML((New-ParaXml $paragraphXml[$idx]))\n    }\n}\n"}
